# Vorstellung_28ster.pptx - "Change Structure and adding new Presentations ..."
#
# Slide 12 ("Mitarbeiter" / "Zuordnung von Ressourcen"): add a new bullet
# at the same outline level (lvl=2, i.e. 3rd level) right after
# "Verplante Projektleiter", reading
# "Offene / Freie Mitarbeiter oder Projektleiter".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Append a brand-new paragraph after the last existing one. Prefixing the
# inserted text with a carriage-return (Chr 13) makes PowerPoint start a
# fresh paragraph instead of just appending a run to the last one; the new
# paragraph inherits the outline level (lvl="2") of the paragraph it is
# split off from, which matches "Verplante Projektleiter".
$lastParagraph = $tr.Paragraphs($tr.Paragraphs().Count)
$lastParagraph.InsertAfter([char]13 + "Offene / Freie Mitarbeiter oder Projektleiter") | Out-Null
